# Add the new data row for 2020-09-01 (SSA raw/clean data) to the
# "out_vars" log sheet, right after the existing last row (93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

# Column A holds dates stored as plain text (e.g. "2020-08-31"), not real
# Excel dates. Force a text number format before assigning the string so
# Excel doesn't auto-convert "2020-09-01" into a date serial, then restore
# the default ("Normal") cell style so the new row matches the style-less
# cells above it.
$aCell = $ws.Cells.Item($row, 1)
$aCell.NumberFormat = "@"
$aCell.Value = "2020-09-01"
$aCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 606036
$ws.Cells.Item($row, 3).Value = 676958
$ws.Cells.Item($row, 4).Value = 77129
$ws.Cells.Item($row, 5).Value = 65241
$ws.Cells.Item($row, 6).Value = 25.36
